# New run commands in NewFlukes2
# Add "Results2" and "Results3" worksheets, each a copy of the existing
# "Results1" sheet (same layout/data as the other Results sheets), to hold
# the data produced by running the refstep program for the report.

$wb = $excel.ActiveWorkbook

# Source worksheet to duplicate.
$src = $wb.Worksheets.Item("Results1")

# Append "Results2" after the current last sheet.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$src.Copy($null, $lastSheet)
$results2 = $wb.Worksheets.Item($wb.Worksheets.Count)
$results2.Name = "Results2"

# Append "Results3" after the (now) last sheet.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$src.Copy($null, $lastSheet)
$results3 = $wb.Worksheets.Item($wb.Worksheets.Count)
$results3.Name = "Results3"
